$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add column C with the sum of columns A and B for each data row (rows 1-11)
for ($r = 1; $r -le 11; $r++) {
    $a = $ws.Cells.Item($r, 1).Value()
    $b = $ws.Cells.Item($r, 2).Value()
    $ws.Cells.Item($r, 3).Value = $a + $b
}

# Update the active selection on the sheet
$ws.Range("D9").Select()
